# Penalty Reward System (unfinished) — shift forecast weeks forward by one
# week and refresh the dependent Summary metrics.
#
# Note: values that look like dates or plain numbers are assigned with a
# leading apostrophe so Excel keeps them as literal Text cells (matching
# the workbook's existing convention of storing these columns as text)
# instead of auto-converting them to date serials / numbers.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet: shift Week_Start_Date forward one week, update MyForecast ---
$ws1.Range("B2").Value = "'2025-01-12"
$ws1.Range("D2").Value = 2
$ws1.Range("B3").Value = "'2025-01-19"
$ws1.Range("D3").Value = 3
$ws1.Range("B4").Value = "'2025-01-26"
$ws1.Range("D4").Value = 3
$ws1.Range("B5").Value = "'2025-02-02"
$ws1.Range("D5").Value = 3
$ws1.Range("B6").Value = "'2025-02-09"
$ws1.Range("D6").Value = 3
$ws1.Range("B7").Value = "'2025-02-16"
$ws1.Range("D7").Value = 3
$ws1.Range("B8").Value = "'2025-02-23"
$ws1.Range("D8").Value = 3
$ws1.Range("B9").Value = "'2025-03-02"
$ws1.Range("D9").Value = 3
$ws1.Range("B10").Value = "'2025-03-09"
$ws1.Range("D10").Value = 3
$ws1.Range("B11").Value = "'2025-03-16"
$ws1.Range("D11").Value = 3
$ws1.Range("B12").Value = "'2025-03-23"
$ws1.Range("D12").Value = 3
$ws1.Range("B13").Value = "'2025-03-30"
$ws1.Range("D13").Value = 3
$ws1.Range("B14").Value = "'2025-04-06"
$ws1.Range("D14").Value = 3
$ws1.Range("B15").Value = "'2025-04-13"
$ws1.Range("D15").Value = 3
$ws1.Range("B16").Value = "'2025-04-20"
$ws1.Range("D16").Value = 3
$ws1.Range("B17").Value = "'2025-04-27"
$ws1.Range("D17").Value = 3

# --- Summary sheet: refresh the metrics that depend on the shifted weeks ---
$ws2.Range("B2").Value  = "'2023-01-15 to 2025-01-05"
$ws2.Range("B8").Value  = "'918 units"
$ws2.Range("B9").Value  = "'49"
$ws2.Range("B10").Value = "'23"
$ws2.Range("B11").Value = "'10"
$ws2.Range("B12").Value = "'3"
$ws2.Range("B13").Value = "'2025-02-09"
$ws2.Range("B14").Value = "'2"
$ws2.Range("B15").Value = "'2025-01-12"
